# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "70.982.99"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.848.82"
$ws.Range("E3").Value = "  +1.36%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'706.20"

# Row 6 - Solana
$ws.Range("D6").Value = "'172.82"
$ws.Range("E6").Value = "  +0.12%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.846.25"
$ws.Range("E7").Value = "  +1.25%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - XRP
$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  -0.70%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.163"
$ws.Range("E10").Value = "  +0.02%  "

# Row 11 - Toncoin
$ws.Range("D11").Value = "'7.36"
$ws.Range("E11").Value = "  -0.44%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.47%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -0.83%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'36.72"
$ws.Range("E14").Value = "  +0.82%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.498.15"
$ws.Range("E15").Value = "  +1.42%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.852.49"
$ws.Range("E16").Value = "  +1.52%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "71.034.71"
$ws.Range("E17").Value = "  +0.36%  "

# Row 18 - Polkadot
$ws.Range("E18").Value = "  -0.22%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +0.72%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'17.37"
$ws.Range("E20").Value = "  -2.96%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'494.50"
$ws.Range("E21").Value = "  +2.49%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'10.64"
$ws.Range("E22").Value = "  -4.09%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  +0.32%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'85.45"
$ws.Range("E24").Value = "  +1.43%  "

# Row 25 - PEPE
$ws.Range("E25").Value = "  +1.87%  "

# Row 26 - RenderToken
$ws.Range("E26").Value = "  +2.00%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").Value = "'12.18"

# Row 28 - Fetch.AI
$ws.Range("E28").Value = "  -3.26%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value = "'3.20"
$ws.Range("E29").Value = "  +2.23%  "

# Row 30 - Dai
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - NEARProtocol
$ws.Range("D31").Value = "'7.52"
$ws.Range("E31").Value = "  -0.14%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  -1.05%  "

# Row 33 - EthereumClassic
$ws.Range("D33").Value = "'29.45"
$ws.Range("E33").Value = "  -0.18%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  -2.20%  "

# Row 35 - RenzoRestakedETH
$ws.Range("D35").Value = "3.805.73"
$ws.Range("E35").Value = "  +1.61%  "

# Row 36 - Aptos
$ws.Range("D36").Value = "'9.17"
$ws.Range("E36").Value = "  -0.61%  "

# Row 37 - Binance-PegBSC-USD
$ws.Range("E37").Value = "  -0.03%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  +0.46%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +7.53%  "

# Row 40 - was Mantle, now Filecoin
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'6.04"
$ws.Range("E40").Value = "  -0.15%  "

# Row 41 - was Filecoin, now Mantle
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "'1.03"
$ws.Range("E41").Value = "  +5.90%  "

# Row 42 - dogwifhat
$ws.Range("E42").Value = "  -3.39%  "

# Row 44 - FirstDigitalUSD
$ws.Range("E44").Value = "  +0.19%  "

# Row 45 - FLOKI
$ws.Range("D45").Value = "'0.000317"
$ws.Range("E45").Value = "  -3.16%  "

# Row 46 - Monero
$ws.Range("D46").Value = "'163.38"
$ws.Range("E46").Value = "  +0.78%  "

# Row 47 - OKB
$ws.Range("D47").Value = "'48.67"
$ws.Range("E47").Value = "  -0.24%  "

# Row 48 - ONDO
$ws.Range("D48").Value = "'1.40"
$ws.Range("E48").Value = "  +0.42%  "

# Row 49 - Bittensor
$ws.Range("D49").Value = "'415.69"
$ws.Range("E49").Value = "  +1.67%  "

# Row 50 - was Cosmos, now TheGraph
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "'0.299"
$ws.Range("E50").Value = "  -1.25%  "

# Row 51 - was TheGraph, now Cosmos
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'8.62"
$ws.Range("E51").Value = "  +0.33%  "
